$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 65600
$ws.Range("E2").Value = 65.09999999999999
$ws.Range("F2").Value = 3.96
$ws.Range("H2").Value = 46
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 60.2
$ws.Range("N2").Value = 66.04328690552585

# Row 3
$ws.Range("K3").Value = 58.8
$ws.Range("N3").Value = 66.04328690552585

# Row 4
$ws.Range("D4").Value = 31250
$ws.Range("E4").Value = 43.9
$ws.Range("F4").Value = 10.62
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 54.8
$ws.Range("N4").Value = 66.04328690552585

# Row 5
$ws.Range("K5").Value = 52.8
$ws.Range("N5").Value = 66.04328690552585

# Row 6
$ws.Range("K6").Value = 51
$ws.Range("N6").Value = 66.04328690552585

# Row 7
$ws.Range("D7").Value = 61700
$ws.Range("E7").Value = 30.2
$ws.Range("F7").Value = 2.83
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 63
$ws.Range("J7").Value = 56
$ws.Range("K7").Value = 51
$ws.Range("N7").Value = 66.04328690552585
